$d = $word.ActiveDocument

# --- Paragraph 2 clean-up -------------------------------------------------
# The whole second paragraph (everything about the POJO getters) is removed,
# but the hidden "_GoBack" bookmark that sits inside it must survive and end
# up at the end of (the now single) paragraph. A plain Range.Delete() that
# spans completely over a bookmark removes the bookmark too, so the deletion
# is split into two pieces: everything after the bookmark, then everything
# before it.

# 1) "color, and getLinkUrl() to return the link URL." (after the bookmark)
$d.Range(194, 241).Delete()

# 2) "To fill the link tag above a POJO expected with methods getLinkText() "
#    + "to return the link text, getLinkColor() to return link " (before the bookmark)
$d.Range(69, 194).Delete()

# 3) Paragraph 2 now contains only the bookmark; merge it into paragraph 1 by
#    deleting the paragraph mark that separates them.
$d.Range(68, 69).Delete()

# --- Paragraph 1 tag clean-up ---------------------------------------------
# Turn on revision tracking for the remaining two edits and accept them
# immediately afterwards: this finalizes the text while stopping the engine
# from silently coalescing the untouched neighbouring runs ("{{", "/", "}}",
# " test") into the edited one, keeping the run layout close to the source.
$d.TrackRevisions = $true

# 4) Trailing space after the final period: ". " -> "."
$d.Range(67, 68).Delete()

# 5) Collapse the link tag attributes down to just "link:" in front of "/"
#    ("text=linkText url=linkUrl color=linkColor " is removed)
$d.Range(16, 58).Delete()

$d.Revisions.AcceptAll() | Out-Null
$d.TrackRevisions = $false
